$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "In progress"
$ws.Range("G6").Value = "Thomas Kosacz"

$ws.Range("G7").Select()
